$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Jeff Wilson (Instagram)", "https://www.instagram.com/profdumpster/?hl=am-et", "Instagram profile."),
    @("Jeff Wilson (YouTube Interview)", "https://www.youtube.com/watch?v=bhuUdCseF3k", "YouTube interview with Jeff Wilson."),
    @("Jeff Wilson (Initialized Blog)", "https://blog.initialized.com/2022/10/video-interview-with-jeff-wilson-ceo-of-jupe/", "Blog interview with Jeff Wilson, CEO of Jupe."),
    @("Jeff Wilson (Miami Ad School)", "https://miamiadschool.com/event/jeff-wilson-jupe/", "Event at Miami Ad School featuring Jeff Wilson."),
    @("Jeff Wilson (Capital Letter)", "https://www.capitalletter.com/p/jupe", "Article on Jupe in Capital Letter.")
)

$startRow = 8
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][2]
}
